$wb = $excel.ActiveWorkbook

# --- "key" sheet: update parameter key/type rows for the new TMT batch ---
$keyWs = $wb.Worksheets.Item("key")

$keyWs.Range("B1").Value = "typelist"

$keyWs.Range("A2").Value = "testTime_tmt"
$keyWs.Range("B2").Value = "time"

$keyWs.Range("A3").Value = "TMT_A"
$keyWs.Range("B3").Value = "num"

$keyWs.Range("A4").Value = "TMT_B"
$keyWs.Range("B4").ClearFormats()
$keyWs.Range("B4").Value = "num"

# --- "all (2)" sheet: the saved selection now highlights column C ---
$all2Ws = $wb.Worksheets.Item("all (2)")
$all2Ws.Range("C1:C1048576").Select()

# Make "key" the active sheet with B2 selected, matching the saved view state
# (activate/select this last so it ends up as the active tab on save)
$keyWs.Activate()
$keyWs.Range("B2").Select()
